$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "data\input\mu_tr_rho.csv"
$ws.Range("C16").Value = "data\input\mu_tr_rho.csv"
$ws.Range("B17").Value = "data\input\h_k_h_amb_10.csv"
$ws.Range("C17").Value = "data\input\h_k_h_amb_10.csv"

$ws.Range("B17").Select()
